$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "Das virale Reel" paragraph
Replace-Text " Die virale Rolle vermittelte erfolgreich den Reiz des Produkts durch ansprechende visuelle und informative Inhalte." `
             " Das virale Reel vermittelte erfolgreich die Attraktivität des Produkts durch ansprechende Visuals und informative Inhalte."

# 2. Influencer Marketing heading + sentence
Replace-Text "Influencer Marketing:" "Influencer-Marketing:"
Replace-Text " Die Macht des Influencer-Marketings kann nicht überstatiert werden." `
             " Die Macht des Influencer-Marketings kann gar nicht hoch genug eingeschätzt werden."

# 3. Geschmacks- und Geschmackssorten heading + sentence
Replace-Text "Geschmacks- und Geschmackssorten:" "Geschmack und Geschmackssorten:"
Replace-Text " Der Ruf von Contoso Protein Plus für köstliche und vielfältige Aromen war ein wichtiger Verkaufspunkt im Viral-Inhalt." `
             " Der Ruf von Contoso Protein Plus, lecker und abwechslungsreich zu sein, war ein wichtiges Verkaufsargument in den viralen Inhalten."

# 4. Gesundheits- und Fitnesstrends sentence
Replace-Text " Der anhaltende Anstieg des Gesundheits- und Fitnessbewusstseins, kombiniert mit einer Zunahme der Anzahl der Menschen, die Trainingsroutinen und aktive Lebensstile annehmen, schuf einen empfänglichen Markt für ein Produkt wie Contoso Protein Plus." `
             " Das ständig wachsende Bewusstsein für Gesundheit und Fitness, verbunden mit einer steigenden Anzahl von Menschen, die eine Trainingsroutine und einen aktiven Lebensstil verfolgen, hat einen empfänglichen Markt für ein Produkt wie Contoso Protein Plus geschaffen."

# 5. Einfache Verfügbarkeit sentence
Replace-Text " Die Barrierefreiheit des Produkts durch verschiedene Onlinehändler hat den Hype weiter gefördert." `
             " Die Verfügbarkeit des Produkts über verschiedene Onlinehändler hat den Hype weiter angeheizt."

# 6. Positive Bewertungen sentence
Replace-Text " Die Rolle war kein Isolierter Fall." " Das Reel war kein Einzelfall."

# 7. Wort des Mundes heading + sentence
Replace-Text "Wort des Mundes:" "Mundpropaganda:"
Replace-Text " Social Media-Plattformen fördern die rasante Verbreitung von Trends durch Mundwort." `
             " Social Media-Plattformen fördern die schnelle Verbreitung von Trends durch Mundpropaganda."
